# Converted Process from Sequence to Flowchart
# -----------------------------------------------------------------------
# Adds three new "lookup helper" columns (Last Column / Fraud Column /
# Bank Details Changed Column) to the "Payment Types" sheet, and updates
# the corresponding formula-template text on the "Sheet Formulas" sheet.

$wb = $excel.ActiveWorkbook

$paymentTypes = $wb.Worksheets.Item("Payment Types")
$sheetFormulas = $wb.Worksheets.Item("Sheet Formulas")

$pilcrow = [char]182

# -------------------------------------------------------------------
# 1. Headers for the two right-most new columns on "Payment Types".
# -------------------------------------------------------------------
$paymentTypes.Range("F1").Value = "Fraud Column"
$paymentTypes.Range("G1").Value = "Bank Details Changed Column"

# -------------------------------------------------------------------
# 2. "Sheet Formulas" sheet: refresh the Cancellation and MR-FIN
#    formula-template text.
# -------------------------------------------------------------------

# Cancellation row (B4): combined M~/N~ block.
$cancellationFormulas = `
  'M~=IFERROR(VLOOKUP(C{0},''Email Refunds''!E:H,4,FALSE),"-")' + $pilcrow + `
  'N~=IFERROR(VLOOKUP(C{0},Fraud!C:G,5,FALSE),"-")'

$sheetFormulas.Range("B4").Value = $cancellationFormulas

# MR-FIN row (B3): K~/L~/M~/N~/O~ block, now referencing 'Refunds - FIN'
# instead of the old malformed MR-FINC:(E) reference.
$mrFinFormulas = `
  'K~=LOOKUP(9.9E+307,--LEFT(MID(I{0},MIN(FIND({1,2,3,4,5,6,7,8,9,0}, $I{0}&"1023456789")),999),ROW(INDIRECT("1:999"))))' + $pilcrow + `
  'L~=RIGHT(I{0},LEN(I{0})-FIND(" ",I{0},1))' + $pilcrow + `
  'M~=IFERROR(VLOOKUP($C{0},donotprocess!$A:$B,{0},FALSE),"-")' + $pilcrow + `
  'N~=IFERROR(VLOOKUP($C{0},Fraud!$C:$G,5,FALSE),"-")' + $pilcrow + `
  'O~=IFERROR(VLOOKUP(C{0},''Refunds - FIN''!C:E,3,0),"-")'

$sheetFormulas.Range("B3").Value = $mrFinFormulas

# -------------------------------------------------------------------
# 3. "Payment Types" sheet: fill in column E (Last Column) top to
#    bottom, then columns F/G row by row.
# -------------------------------------------------------------------

$paymentTypes.Range("E1").Value = "Last Column"
$paymentTypes.Range("E2").Value = "R"
$paymentTypes.Range("E3").Value = "S"
$paymentTypes.Range("E4").Value = "S"
$paymentTypes.Range("E5").Value = "S"
$paymentTypes.Range("E6").Value = "P"
$paymentTypes.Range("E7").Value = "P"

$paymentTypes.Range("F2").Value = "FRAUD Check"
$paymentTypes.Range("G2").Value = "Do not process"

$paymentTypes.Range("F3").Value = "refunded by Fraud team?"
$paymentTypes.Range("G3").Value = "Manual refunded Yes/No"

$paymentTypes.Range("F4").Value = "refunded by Fraud team?"
$paymentTypes.Range("G4").Value = "Manual refunded Yes/No"

$paymentTypes.Range("F5").Value = "refunded by Fraud team?"
$paymentTypes.Range("G5").Value = "Manual refunded Yes/No"

$paymentTypes.Range("F6").Value = "Fraud"
$paymentTypes.Range("G6").Value = "Manual refund"

$paymentTypes.Range("F7").Value = "Fraud"
$paymentTypes.Range("G7").Value = "Manual refund"

# Column widths for the new columns (best effort match of the authored
# widths: E=28.7109375, F=23.7109375 (bestFit), G=27.85546875 (bestFit))
$paymentTypes.Columns.Item(5).ColumnWidth = 27.83
$paymentTypes.Columns.Item(6).ColumnWidth = 22.83
$paymentTypes.Columns.Item(7).ColumnWidth = 26.92

# -------------------------------------------------------------------
# 4. View state: make "Payment Types" the active sheet/tab, with F4
#    selected; leave "Sheet Formulas" selection at C8.
# -------------------------------------------------------------------

$sheetFormulas.Range("C8").Select()
$paymentTypes.Activate()
$paymentTypes.Range("F4").Select()
